# Update the "Nama Sub Kategori" (column C) labels for a handful of
# sub-categories on the "Sub Kategori Produk" sheet: the old entries
# duplicated the parent word ("Minuman"/"Perabotan") in the label, the
# new entries drop it since it's implied by context.
#
# Write order matches the order the new shared-string entries were
# appended in the authored workbook (rows 56/55 first, then 49/50/51).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C56").Value = "Dapur"
$ws.Range("C55").Value = " Masak"
$ws.Range("C49").Value = "Non Alkohol"
$ws.Range("C50").Value = "Fermentasi"
$ws.Range("C51").Value = "Alkohol"

# Restore the cursor/viewport position recorded in the saved file.
[void]$ws.Range("C30").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 19
$win.ScrollColumn = 1
